$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabla#2 - Tokens")

# New Lexema/Token text values for rows 3-10 (column B and E)
$ws.Range("B3").Value = "fin"
$ws.Range("E3").Value = "fin"

$ws.Range("B4").Value = "fin"
$ws.Range("E4").Value = "fin"

$ws.Range("B5").Value = "inicio"
$ws.Range("E5").Value = "inicio"

$ws.Range("B6").Value = "inicio"
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = "inicio"

$ws.Range("B7").Value = "peso"
$ws.Range("C7").Value = 8
$ws.Range("E7").Value = "peso"

$ws.Range("B8").Value = "peso"
$ws.Range("D8").Value = 15
$ws.Range("E8").Value = "peso"

$ws.Range("B9").Value = "nombre"
$ws.Range("E9").Value = "nombre"

$ws.Range("B10").Value = "nombre"
$ws.Range("E10").Value = "nombre"
